# Dependency Injection & Providers - 1 - Overview
#
# The original Sheet1 (module-review sequencing list) is preserved as a new
# "Sheet2", and Sheet1 is repurposed to hold the new "Module Category"
# tracking grid (category / expertise level / usage frequency / review
# frequency, plus the full module list rows 4-17).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Clone the old Sheet1 data onto a brand-new "Sheet2" ----------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("B6").Value  = "Unit Testing - 1 - Overview"
$ws2.Range("B7").Value  = "Unit Testing - 2 -Testing Classes & Pipes"
$ws2.Range("B8").Value  = "Unit Testing - 3 -Jasmine & Karma"
$ws2.Range("B9").Value  = "Unit Testing - 4 - Testing with Mocks & Spies"
$ws2.Range("B13").Value = "Quickstart - Nesting Components & Inputs"
$ws2.Range("B14").Value = "Quickstart - User Interaction & Outputs"
$ws2.Range("B15").Value = "Unit Testing - 9 - Testing Components - todo"

# --- 2. Wipe Sheet1 and write the new tracking-grid content -----------------
$ws1.Cells.Clear()

$ws1.Range("A1").Value = "Module Category"
$ws1.Range("B1").Value = "Your level of expertise"
$ws1.Range("C1").Value = "Level of use in daily life"
$ws1.Range("D1").Value = "Frequency to review so that we will not forget the contents"

$ws1.Range("C2").Value = "Basic/Intermediate/Expert"

$ws1.Range("A4").Value  = "1.quickstart"
$ws1.Range("A5").Value  = "2.es6-typescript"
$ws1.Range("A6").Value  = "3.angular-cli"
$ws1.Range("A7").Value  = "4.components"
$ws1.Range("A8").Value  = "5.built-in-directives"
$ws1.Range("A9").Value  = "6.custom-directives"
$ws1.Range("A10").Value = "7.reactive-programming-with-rxjs"
$ws1.Range("A11").Value = "8.pipes"
$ws1.Range("A12").Value = "9.forms"
$ws1.Range("A13").Value = "10.dependency-injection-and-providers"
$ws1.Range("A14").Value = "11.HTTP"
$ws1.Range("A15").Value = "12.routing"
$ws1.Range("A16").Value = "13.unit-testing"
$ws1.Range("A17").Value = "14.advanced-topics"

# --- 3. Column widths (best-fit, as in the authored workbook) --------------
$ws1.Columns.Item(1).ColumnWidth = 32.666666666666664
$ws1.Columns.Item(2).ColumnWidth = 18.5
$ws1.Columns.Item(3).ColumnWidth = 22
$ws1.Columns.Item(4).ColumnWidth = 49.333333333333336

# --- 4. Selections / active sheet -------------------------------------------
# Select Sheet2's range first, then Sheet1's - last selection wins the
# "active sheet" (tabSelected) state, and Sheet1 should stay active.
$ws2.Range("B5:G19").Select()
$ws1.Range("B2").Select()
